# Apply "current state (new derivation)" edits:
#  1. gens!M2:M34  -> overwrite genCostFactor column with constant 1
#  2. splice_lines -> delete the two data rows (2 and 3), leaving only the header
#  3. Make "gens" the active/selected sheet (was "busses")

$wb = $excel.ActiveWorkbook

# --- 1. gens sheet: set genCostFactor column (M) to 1 for all data rows ---
$gens = $wb.Worksheets.Item("gens")
$gens.Range("M2:M34").Value = 1

# --- 2. splice_lines sheet: remove the two data rows entirely ---
$splice = $wb.Worksheets.Item("splice_lines")
$splice.Rows.Item(2).Resize(2).Delete()
$splice.Range("A2").Select()
$splice.Range("A2:XFD3").Select()

# --- 3. Switch the active sheet to gens, select M2:M34 there ---
$gens.Select()
$gens.Range("M2:M34").Select()
